$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") should carry the same style as
# the existing header cells (e.g. H1): bold font, thin border, centered.
# Copy the formatting from H1 (paste formats only) before writing the text,
# so the new cells reuse the existing style entry instead of minting a new one.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data cells for rows 2 and 3
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9
